# "login account, fixed bug"
# - Updates the saved login password in the Config sheet and moves the
#   cursor to the next field.
# - Corrects the Order/Result tracking numbers on the TrackProduct sheet
#   (a second run against item #1 instead of a new item #2, a later
#   timestamp, and the resulting pass/fail statuses), and moves the
#   cursor off the stale J3 cell.

$wb = $excel.ActiveWorkbook

# ---- Config sheet ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Config")

# Updated account password.
$ws2.Range("E7").Value = "qyqqyq123"

# Move selection to the next field.
$ws2.Range("E8").Select() | Out-Null

# ---- TrackProduct sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item("TrackProduct")
$ws1.Activate()

# Row 2: this run actually failed (price target not met) at a later time.
$ws1.Range("F2").Value = "fail"
$ws1.Range("H2").Value = "25-May-2021 17h01m"
$ws1.Range("K2").Value = 0
$ws1.Range("L2").Value = "fail"

# Row 3: it was a second attempt against item #1 (not a brand-new item #2),
# recorded later, and order-processed count bumped to 1.
$ws1.Range("C3").Value = 1
$ws1.Range("H3").Value = "25-May-2021 17h07m"
$ws1.Range("K3").Value = 1

# Result table ref spans H1:L3, keep its AutoFilter range in sync.
$ws1.ListObjects("Result").AutoFilter.Range = $ws1.Range("H1:L3")

# Move selection off the old J3 cell; TrackProduct stays the active tab.
$ws1.Range("F3").Select() | Out-Null
